# A new weekly price observation (Membrillo, Primera, Vega Modelo de
# Temuco) was added to the dataset. In the source sheet this shows up as
# a brand-new row inserted right before the former row 245, which pushes
# every row from the old 245 down through the old 312 down by one
# (old row 245 -> new row 246, ..., old row 312 -> new row 313).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 245, shifting rows 245:312 down to 246:313.
$ws.Rows.Item(245).Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A245").Value = 10
$ws.Range("B245").Value = "Vega Modelo de Temuco"
$ws.Range("C245").Value = "La Araucanía"
$ws.Range("D245").Value = 45135
$ws.Range("E245").Value = 9
$ws.Range("F245").Value = "Fruta"
$ws.Range("G245").Value = 100104
$ws.Range("H245").Value = "Frutos de pepita"
$ws.Range("I245").Value = 100104003
$ws.Range("J245").Value = "Membrillo"
$ws.Range("K245").Value = "Champion"
$ws.Range("L245").Value = "Primera"
$ws.Range("M245").Value = 50
$ws.Range("N245").Value = 16000
$ws.Range("O245").Value = 16000
$ws.Range("P245").Value = 16000
$ws.Range("Q245").Value = "$/bandeja 18 kilos granel"
$ws.Range("R245").Value = "Región de O'Higgins"
$ws.Range("S245").Value = 889
$ws.Range("T245").Value = 18
